$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - daily spot price update
$ws.Range("A2").Value = 46047
$ws.Range("B2").Value = 21.53
$ws.Range("C2").Value = 16.03
$ws.Range("D2").Value = 5.46
$ws.Range("E2").Value = 3.78
$ws.Range("F2").Value = 3.78
$ws.Range("G2").Value = 3.78
$ws.Range("H2").Value = 3.78
$ws.Range("I2").Value = 3.78
$ws.Range("J2").Value = 3.96
$ws.Range("K2").Value = 3.78
$ws.Range("L2").Value = 1.45
$ws.Range("M2").Value = 0.32
$ws.Range("N2").Value = 0.04
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0.57
$ws.Range("S2").Value = 3.53
$ws.Range("T2").Value = 28.92
$ws.Range("U2").Value = 35.75
$ws.Range("V2").Value = 36.44
$ws.Range("W2").Value = 37.35
$ws.Range("X2").Value = 34.49
$ws.Range("Y2").Value = 21.56
$ws.Range("Z2").Value = 11.25
$ws.Range("AB2").Value = 32.46
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 36.89
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 32.34
$ws.Range("AG2").Value = "2h-17h"
